# Update image path strings: remove the "Examples\" path segment
# from "..\CodeSnippets\Examples\Snippets\*.png" -> "..\CodeSnippets\Snippets\*.png"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = @("A2", "D2", "G2", "J2", "M2", "P2", "S2", "A3", "D3", "G3", "J3", "M3", "P3", "S3")

foreach ($cellAddr in $cells) {
    $range = $ws.Range($cellAddr)
    $oldValue = [string]$range.Value2
    $newValue = $oldValue.Replace("..\CodeSnippets\Examples\Snippets\", "..\CodeSnippets\Snippets\")
    $range.Value2 = $newValue
}
